$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 122.253015
$ws.Cells.Item(2, 8).Value = 366.759045
$ws.Cells.Item(2, 9).Value = 0.1988639364328829
$ws.Cells.Item(2, 10).Value = 0.1988639364328829
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 136.6884486666667
$ws.Cells.Item(2, 14).Value = 410.065346
$ws.Cells.Item(2, 15).Value = 0.7423691870207686
$ws.Cells.Item(2, 16).Value = 0.7423691870207685
$ws.Cells.Item(2, 17).Value = 16710.57496517273
$ws.Cells.Item(2, 18).Value = 150395.1746865546
$ws.Cells.Item(2, 19).Value = 0.1476304588174291
$ws.Cells.Item(2, 20).Value = 0.147630458817429

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 122.253015
$ws.Cells.Item(3, 8).Value = 366.759045
$ws.Cells.Item(3, 9).Value = 0.1988639364328829
$ws.Cells.Item(3, 10).Value = 0.1988639364328829
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 0.8952453333333334
$ws.Cells.Item(3, 14).Value = 2.685736
$ws.Cells.Item(3, 15).Value = 0.004862170554817893
$ws.Cells.Item(3, 16).Value = 0.004862170554817893
$ws.Cells.Item(3, 17).Value = 109.44644116468
$ws.Cells.Item(3, 18).Value = 985.0179704821202
$ws.Cells.Item(3, 19).Value = 0.0009669103761391404
$ws.Cells.Item(3, 20).Value = 0.0009669103761391402

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 122.253015
$ws.Cells.Item(4, 8).Value = 366.759045
$ws.Cells.Item(4, 9).Value = 0.1988639364328829
$ws.Cells.Item(4, 10).Value = 0.1988639364328829
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 41.63761133333333
$ws.Cells.Item(4, 14).Value = 124.912834
$ws.Cells.Item(4, 15).Value = 0.2261381995079395
$ws.Cells.Item(4, 16).Value = 0.2261381995079395
$ws.Cells.Item(4, 17).Value = 5090.32352289817
$ws.Cells.Item(4, 18).Value = 45812.91170608353
$ws.Cells.Item(4, 19).Value = 0.04497073253199346
$ws.Cells.Item(4, 20).Value = 0.04497073253199346

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 122.253015
$ws.Cells.Item(5, 8).Value = 366.759045
$ws.Cells.Item(5, 9).Value = 0.1988639364328829
$ws.Cells.Item(5, 10).Value = 0.1988639364328829
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 4.903320333333333
$ws.Cells.Item(5, 14).Value = 14.709961
$ws.Cells.Item(5, 15).Value = 0.02663044291647413
$ws.Cells.Item(5, 16).Value = 0.02663044291647413
$ws.Cells.Item(5, 17).Value = 599.445694260805
$ws.Cells.Item(5, 18).Value = 5395.011248347245
$ws.Cells.Item(5, 19).Value = 0.005295834707321228
$ws.Cells.Item(5, 20).Value = 0.005295834707321226

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 132.5447616666667
$ws.Cells.Item(6, 8).Value = 397.634285
$ws.Cells.Item(6, 9).Value = 0.2156050961899926
$ws.Cells.Item(6, 10).Value = 0.2156050961899926
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 136.6884486666667
$ws.Cells.Item(6, 14).Value = 410.065346
$ws.Cells.Item(6, 15).Value = 0.7423691870207686
$ws.Cells.Item(6, 16).Value = 0.7423691870207685
$ws.Cells.Item(6, 17).Value = 18117.33785110973
$ws.Cells.Item(6, 18).Value = 163056.0406599876
$ws.Cells.Item(6, 19).Value = 0.1600585799760994
$ws.Cells.Item(6, 20).Value = 0.1600585799760994

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 132.5447616666667
$ws.Cells.Item(7, 8).Value = 397.634285
$ws.Cells.Item(7, 9).Value = 0.2156050961899926
$ws.Cells.Item(7, 10).Value = 0.2156050961899926
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 0.8952453333333334
$ws.Cells.Item(7, 14).Value = 2.685736
$ws.Cells.Item(7, 15).Value = 0.004862170554817893
$ws.Cells.Item(7, 16).Value = 0.004862170554817893
$ws.Cells.Item(7, 17).Value = 118.6600793398622
$ws.Cells.Item(7, 18).Value = 1067.94071405876
$ws.Cells.Item(7, 19).Value = 0.001048308750163661
$ws.Cells.Item(7, 20).Value = 0.001048308750163661

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 132.5447616666667
$ws.Cells.Item(8, 8).Value = 397.634285
$ws.Cells.Item(8, 9).Value = 0.2156050961899926
$ws.Cells.Item(8, 10).Value = 0.2156050961899926
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 41.63761133333333
$ws.Cells.Item(8, 14).Value = 124.912834
$ws.Cells.Item(8, 15).Value = 0.2261381995079395
$ws.Cells.Item(8, 16).Value = 0.2261381995079395
$ws.Cells.Item(8, 17).Value = 5518.847270545965
$ws.Cells.Item(8, 18).Value = 49669.62543491369
$ws.Cells.Item(8, 19).Value = 0.04875654825714102
$ws.Cells.Item(8, 20).Value = 0.04875654825714102

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 132.5447616666667
$ws.Cells.Item(9, 8).Value = 397.634285
$ws.Cells.Item(9, 9).Value = 0.2156050961899926
$ws.Cells.Item(9, 10).Value = 0.2156050961899926
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 4.903320333333333
$ws.Cells.Item(9, 14).Value = 14.709961
$ws.Cells.Item(9, 15).Value = 0.02663044291647413
$ws.Cells.Item(9, 16).Value = 0.02663044291647413
$ws.Cells.Item(9, 17).Value = 649.9094249569872
$ws.Cells.Item(9, 18).Value = 5849.184824612885
$ws.Cells.Item(9, 19).Value = 0.005741659206588512
$ws.Cells.Item(9, 20).Value = 0.005741659206588511

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 320.0894206666666
$ws.Cells.Item(10, 8).Value = 960.2682619999999
$ws.Cells.Item(10, 9).Value = 0.5206762565675317
$ws.Cells.Item(10, 10).Value = 0.5206762565675317
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 136.6884486666667
$ws.Cells.Item(10, 14).Value = 410.065346
$ws.Cells.Item(10, 15).Value = 0.7423691870207686
$ws.Cells.Item(10, 16).Value = 0.7423691870207685
$ws.Cells.Item(10, 17).Value = 43752.52634553873
$ws.Cells.Item(10, 18).Value = 393772.7371098486
$ws.Cells.Item(10, 19).Value = 0.3865340092890556
$ws.Cells.Item(10, 20).Value = 0.3865340092890556

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 320.0894206666666
$ws.Cells.Item(11, 8).Value = 960.2682619999999
$ws.Cells.Item(11, 9).Value = 0.5206762565675317
$ws.Cells.Item(11, 10).Value = 0.5206762565675317
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 0.8952453333333334
$ws.Cells.Item(11, 14).Value = 2.685736
$ws.Cells.Item(11, 15).Value = 0.004862170554817893
$ws.Cells.Item(11, 16).Value = 0.004862170554817893
$ws.Cells.Item(11, 17).Value = 286.5585601012036
$ws.Cells.Item(11, 18).Value = 2579.027040910832
$ws.Cells.Item(11, 19).Value = 0.002531616763275459
$ws.Cells.Item(11, 20).Value = 0.002531616763275459

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 320.0894206666666
$ws.Cells.Item(12, 8).Value = 960.2682619999999
$ws.Cells.Item(12, 9).Value = 0.5206762565675317
$ws.Cells.Item(12, 10).Value = 0.5206762565675317
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 41.63761133333333
$ws.Cells.Item(12, 14).Value = 124.912834
$ws.Cells.Item(12, 15).Value = 0.2261381995079395
$ws.Cells.Item(12, 16).Value = 0.2261381995079395
$ws.Cells.Item(12, 17).Value = 13327.7588896305
$ws.Cells.Item(12, 18).Value = 119949.8300066745
$ws.Cells.Item(12, 19).Value = 0.1177447911867155
$ws.Cells.Item(12, 20).Value = 0.1177447911867155

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 320.0894206666666
$ws.Cells.Item(13, 8).Value = 960.2682619999999
$ws.Cells.Item(13, 9).Value = 0.5206762565675317
$ws.Cells.Item(13, 10).Value = 0.5206762565675317
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 4.903320333333333
$ws.Cells.Item(13, 14).Value = 14.709961
$ws.Cells.Item(13, 15).Value = 0.02663044291647413
$ws.Cells.Item(13, 16).Value = 0.02663044291647413
$ws.Cells.Item(13, 17).Value = 1569.500964839753
$ws.Cells.Item(13, 18).Value = 14125.50868355778
$ws.Cells.Item(13, 19).Value = 0.01386583932848509
$ws.Cells.Item(13, 20).Value = 0.01386583932848509

$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 39.86989333333333
$ws.Cells.Item(14, 8).Value = 119.60968
$ws.Cells.Item(14, 9).Value = 0.06485471080959287
$ws.Cells.Item(14, 10).Value = 0.06485471080959287
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 136.6884486666667
$ws.Cells.Item(14, 14).Value = 410.065346
$ws.Cells.Item(14, 15).Value = 0.7423691870207686
$ws.Cells.Item(14, 16).Value = 0.7423691870207685
$ws.Cells.Item(14, 17).Value = 5449.753868238809
$ws.Cells.Item(14, 18).Value = 49047.78481414928
$ws.Cells.Item(14, 19).Value = 0.04814613893818451
$ws.Cells.Item(14, 20).Value = 0.04814613893818451

$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 39.86989333333333
$ws.Cells.Item(15, 8).Value = 119.60968
$ws.Cells.Item(15, 9).Value = 0.06485471080959287
$ws.Cells.Item(15, 10).Value = 0.06485471080959287
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 0.8952453333333334
$ws.Cells.Item(15, 14).Value = 2.685736
$ws.Cells.Item(15, 15).Value = 0.004862170554817893
$ws.Cells.Item(15, 16).Value = 0.004862170554817893
$ws.Cells.Item(15, 17).Value = 35.69333594716445
$ws.Cells.Item(15, 18).Value = 321.24002352448
$ws.Cells.Item(15, 19).Value = 0.0003153346652396322
$ws.Cells.Item(15, 20).Value = 0.0003153346652396322

$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 39.86989333333333
$ws.Cells.Item(16, 8).Value = 119.60968
$ws.Cells.Item(16, 9).Value = 0.06485471080959287
$ws.Cells.Item(16, 10).Value = 0.06485471080959287
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 41.63761133333333
$ws.Cells.Item(16, 14).Value = 124.912834
$ws.Cells.Item(16, 15).Value = 0.2261381995079395
$ws.Cells.Item(16, 16).Value = 0.2261381995079395
$ws.Cells.Item(16, 17).Value = 1660.087122514791
$ws.Cells.Item(16, 18).Value = 14940.78410263312
$ws.Cells.Item(16, 19).Value = 0.01466612753208943
$ws.Cells.Item(16, 20).Value = 0.01466612753208943

$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 39.86989333333333
$ws.Cells.Item(17, 8).Value = 119.60968
$ws.Cells.Item(17, 9).Value = 0.06485471080959287
$ws.Cells.Item(17, 10).Value = 0.06485471080959287
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 4.903320333333333
$ws.Cells.Item(17, 14).Value = 14.709961
$ws.Cells.Item(17, 15).Value = 0.02663044291647413
$ws.Cells.Item(17, 16).Value = 0.02663044291647413
$ws.Cells.Item(17, 17).Value = 195.4948586691644
$ws.Cells.Item(17, 18).Value = 1759.45372802248
$ws.Cells.Item(17, 19).Value = 0.001727109674079301
$ws.Cells.Item(17, 20).Value = 0.001727109674079301
